$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the entire row 436 ("「ヨユウなきユウヨ」"), shifting all subsequent
# rows up by one. This reduces the used range from A1:C610 to A1:C609.
$ws.Rows.Item(436).Delete()
